$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.900809049606323
$ws.Range("B1").Value = 4.947888374328613
$ws.Range("C1").Value = 3.667495250701904
$ws.Range("D1").Value = 0.4506295025348663
$ws.Range("E1").Value = 0.2405924648046494
